$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.533.37"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = "'2.457.75"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("D5").Value = "'314.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").Value = "'91.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("D7").Value = "'0.547"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.47%  '
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("E9").Value = '  +4.70%  '
$ws.Range("D10").Value = "'32.40"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.28%  '
$ws.Range("D11").Value = "'0.0795"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.25%  '
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("D13").Value = "'2.838.01"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = "'15.85"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.74%  '
$ws.Range("D16").Value = "'2.469.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.24%  '
$ws.Range("D17").Value = "'0.768"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.45%  '
$ws.Range("D18").Value = "'41.523.02"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").Value = "'6.45"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.36%  '
$ws.Range("E20").Value = '  +3.18%  '
$ws.Range("D21").Value = "'70.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("D22").Value = "'11.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.08%  '
$ws.Range("D23").Value = "'236.63"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.07%  '
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("E26").Value = '  +2.43%  '
$ws.Range("D27").Value = "'24.22"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.12%  '
$ws.Range("D28").Value = "'2.26"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.69%  '
$ws.Range("E29").Value = '  +1.42%  '
$ws.Range("D30").Value = "'34.94"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("D31").Value = "'155.47"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("E32").Value = '  +2.68%  '
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("D34").Value = "'0.0758"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("D35").Value = "'17.41"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.19%  '
$ws.Range("D36").Value = "'2.41"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.11%  '
$ws.Range("D37").Value = "'2.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("E38").Value = '  +2.81%  '
$ws.Range("E39").Value = '  +2.11%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("D43").Value = "'1.967.20"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.09%  '
$ws.Range("D44").Value = "'0.0281"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.91%  '
$ws.Range("D45").Value = "'18.35"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -9.45%  '
$ws.Range("D46").Value = "'2.89"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("D47").Value = "'8.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +4.11%  '
$ws.Range("D48").Value = "'2.697.65"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("D49").Value = "'96.23"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("D50").Value = "'66.09"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").Value = "'71.69"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.52%  '
